$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 2, shifting existing rows (and their formatting) down.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with its data.
$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "ETH"
$ws.Range("T2").Value = 748.94669999999996

# Update the view's selection as recorded in the saved workbook.
$ws.Range("Q2:S10").Select()
